$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
Write-Host ("Dimension before: " + $ws.UsedRange.Address())
